$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "29.488.38"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "1.918.12"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("E4").Value = "  +0.72%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "325.43"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4801"
$ws.Range("E7").Value = "  -0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4052"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08211"
$ws.Range("E9").Value = "  +0.75%  "
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.43"
$ws.Range("E11").Value = "  -0.24%  "
$ws.Range("D12").Value = "1.931.72"
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.051"
$ws.Range("E13").Value = "  +0.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.224"
$ws.Range("E14").Value = "  +1.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.33"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06854"
$ws.Range("E16").Value = "  +1.53%  "
$ws.Range("E17").Value = "  +0.65%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001037"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.53"
$ws.Range("E19").Value = "  -0.88%  "
$ws.Range("E20").Value = "  +0.57%  "
$ws.Range("D21").Value = "29.492.74"
$ws.Range("E21").Value = "  -0.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.668"
$ws.Range("E22").Value = "  +1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.89"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.192"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "2.152.88"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.15"
$ws.Range("E26").Value = "  +1.16%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.505"
$ws.Range("E27").Value = "  +2.77%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("E29").Value = "  -0.65%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.41"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.013"
$ws.Range("E31").Value = "  -2.41%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09603"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.615"
$ws.Range("E33").Value = "  +1.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.557"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.369"
$ws.Range("E35").Value = "  -2.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06314"
$ws.Range("E36").Value = "  +3.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02278"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("E38").Value = "  +0.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5928"
$ws.Range("E39").Value = "  -0.58%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.71"
$ws.Range("E40").Value = "  +1.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.873"
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.384"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.40"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07477"
$ws.Range("E46").Value = "  -3.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5555"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "117.90"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.424"
$ws.Range("E50").Value = "  +3.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.91"
$ws.Range("E51").Value = "  -0.99%  "
